# Update the "Förändrad" (Changed) date column (C) from 2024-12-14 (45640)
# to 2024-12-15 (45641) for all data rows (2 through 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 34
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45640) {
        $cell.Value2 = 45641
    }
}
